# Updated cryptos list (Price + Volume(1h) columns) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as plain text, preserving the original cell
# style (Price cells hold values like "29.258.06" that Excel would otherwise
# auto-convert/mis-parse as a number).
function Set-TextValue($range, $text) {
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $savedStyle
}

Set-TextValue $ws.Range("D2") '29.258.06'
Set-TextValue $ws.Range("E2") '  +0.14%  '
Set-TextValue $ws.Range("D3") '1.899.54'
Set-TextValue $ws.Range("E3") '  -0.33%  '
Set-TextValue $ws.Range("E4") '  -0.20%  '
Set-TextValue $ws.Range("D5") '326.25'
Set-TextValue $ws.Range("E5") '  -0.42%  '
Set-TextValue $ws.Range("D6") '1.001'
Set-TextValue $ws.Range("E6") '  -0.16%  '
Set-TextValue $ws.Range("E7") '  -0.51%  '
Set-TextValue $ws.Range("D8") '0.3918'
Set-TextValue $ws.Range("E8") '  -0.95%  '
Set-TextValue $ws.Range("D9") '0.07876'
Set-TextValue $ws.Range("E9") '  -1.16%  '
Set-TextValue $ws.Range("D10") '0.9879'
Set-TextValue $ws.Range("E10") '  -1.58%  '
Set-TextValue $ws.Range("D11") '21.84'
Set-TextValue $ws.Range("E11") '  -2.25%  '
Set-TextValue $ws.Range("D12") '1.909.14'
Set-TextValue $ws.Range("E12") '  +1.61%  '
Set-TextValue $ws.Range("D13") '7.061'
Set-TextValue $ws.Range("E13") '  -1.06%  '
Set-TextValue $ws.Range("E14") '  -0.97%  '
Set-TextValue $ws.Range("D15") '0.06981'
Set-TextValue $ws.Range("E15") '  +0.47%  '
Set-TextValue $ws.Range("E16") '  -0.67%  '
Set-TextValue $ws.Range("E17") '  -0.17%  '
Set-TextValue $ws.Range("D18") '0.000009972'
Set-TextValue $ws.Range("E18") '  -1.53%  '
Set-TextValue $ws.Range("E19") '  -0.88%  '
Set-TextValue $ws.Range("D20") '1.001'
Set-TextValue $ws.Range("E20") '  -0.20%  '
Set-TextValue $ws.Range("D21") '29.271.51'
Set-TextValue $ws.Range("E21") '  +0.12%  '
Set-TextValue $ws.Range("D22") '5.296'
Set-TextValue $ws.Range("D23") '11.08'
Set-TextValue $ws.Range("E23") '  -0.08%  '
Set-TextValue $ws.Range("D24") '2.099'
Set-TextValue $ws.Range("E24") '  +1.97%  '
Set-TextValue $ws.Range("D25") '156.09'
Set-TextValue $ws.Range("E25") '  -0.40%  '
Set-TextValue $ws.Range("E26") '  -0.79%  '
Set-TextValue $ws.Range("D27") '6.025'
Set-TextValue $ws.Range("E27") '  +2.20%  '
Set-TextValue $ws.Range("D28") '118.46'
Set-TextValue $ws.Range("E28") '  -0.96%  '
Set-TextValue $ws.Range("D29") '1.888'
Set-TextValue $ws.Range("E29") '  -5.78%  '
Set-TextValue $ws.Range("D30") '0.09358'
Set-TextValue $ws.Range("E30") '  -0.85%  '
Set-TextValue $ws.Range("D31") '0.9034'
Set-TextValue $ws.Range("E31") '  -2.39%  '
Set-TextValue $ws.Range("D32") '5.256'
Set-TextValue $ws.Range("E32") '  -1.86%  '
Set-TextValue $ws.Range("D33") '1.322'
Set-TextValue $ws.Range("E33") '  -1.85%  '
Set-TextValue $ws.Range("D34") '3.213'
Set-TextValue $ws.Range("E34") '  -1.51%  '
Set-TextValue $ws.Range("D35") '1.189'
Set-TextValue $ws.Range("E35") '  +1.28%  '
Set-TextValue $ws.Range("D36") '0.05773'
Set-TextValue $ws.Range("E36") '  -1.52%  '
Set-TextValue $ws.Range("D37") '0.02084'
Set-TextValue $ws.Range("E37") '  -1.12%  '
Set-TextValue $ws.Range("E38") '  -0.18%  '
Set-TextValue $ws.Range("D39") '7.716'
Set-TextValue $ws.Range("E40") '  -0.93%  '
Set-TextValue $ws.Range("D41") '0.1788'
Set-TextValue $ws.Range("E41") '  -1.43%  '
Set-TextValue $ws.Range("D42") '9.711'
Set-TextValue $ws.Range("E42") '  -3.23%  '
Set-TextValue $ws.Range("D43") '11.89'
Set-TextValue $ws.Range("E43") '  -1.01%  '
Set-TextValue $ws.Range("D44") '0.5355'
Set-TextValue $ws.Range("E44") '  -1.50%  '
Set-TextValue $ws.Range("E45") '  -2.12%  '
Set-TextValue $ws.Range("D46") '0.07026'
Set-TextValue $ws.Range("E46") '  -1.01%  '
Set-TextValue $ws.Range("E47") '  -1.99%  '
Set-TextValue $ws.Range("D48") '2.569'
Set-TextValue $ws.Range("E48") '  -0.43%  '
Set-TextValue $ws.Range("D49") '112.94'
Set-TextValue $ws.Range("E49") '  +0.77%  '
Set-TextValue $ws.Range("D50") '1.057'
Set-TextValue $ws.Range("E50") '  -1.09%  '
Set-TextValue $ws.Range("D51") '0.2901'
Set-TextValue $ws.Range("E51") '  +0.07%  '
